$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------
# Alternative scenario step renumbered: "3.1 Requisita a inserção de um
# novo valor" / "3.2 Regressa a 1" collapse into a single "3.1 Regressa a 1"
$ws.Range("D13").Value = "3.1 Regressa a 1"
$ws.Range("D14").Value = ""

# Exception label now references step 1 instead of step 3.1
$ws.Range("B16").Value = "Exceção 2 [não insere novo valor] (passo 1)"

# --- Formatting updates ------------------------------------------------
# Shrink font used by the "Cenário Alternativo 1" merged label (B13:B15)
$ws.Range("B13:B15").Font.Size = 10

# Row 13 no longer needs the extra wrapped line now that the text is shorter
$ws.Rows.Item(13).RowHeight = 19.5

# --- View state ----------------------------------------------------------
$excel.ActiveWindow.Zoom = 90
$ws.Range("B18").Select()
